$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.441.02"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").Value = "1.805.19"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4478"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3757"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07523"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.625"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.306"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").Value = "1.802.20"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001092"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06778"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9995"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.318"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.65%  "

$ws.Range("D23").Value = "28.429.81"
$ws.Range("E23").Value = "  +1.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.407"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.46%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.19%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.57%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.357"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.81%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "2.006.46"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.260"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.90%  "

$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.008"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.80%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.825"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.92%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09334"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.05%  "

$ws.Range("B35").Value = "Algorand"
$ws.Range("C35").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2263"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.69%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.79%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06359"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02333"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.63%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6588"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.158"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.93%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.211"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.02%  "

$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.450"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.94%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.105"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9986"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6077"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.814"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.035"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07092"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.157"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.13%  "
